# Append a new row (row 4) to the "List To Bring Down" sheet, mirroring the
# existing "Oil" row (row 2): ProductName, Quantity, Big StockRoom = Y,
# Small Stockroom / Cooler left blank (stored as empty text, matching the
# other blank cells already present on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Oil"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Y"

# Writing a bare apostrophe stores an empty string (Excel's "quote prefix"
# convention) instead of clearing the cell back to blank/null. Re-apply the
# neighbouring (unstyled) cell's style afterwards so we don't pick up the
# quote-prefix formatting flag that the apostrophe entry would otherwise add.
$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = $ws.Range("D2").Style

$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = $ws.Range("E2").Style
